$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add sample data in I2 first (new cell -> new shared string entry)
$ws.Range("I2").Value = "2022, 2023 - 2nd"

# Rename header I1 from AYCODE to BATCHYEAR (reuses/overwrites its shared string slot)
$ws.Range("I1").Value = "BATCHYEAR"

# Size the new column I to fit its header/content (bestfit ~15 chars wide)
$ws.Columns("I:I").ColumnWidth = 14.166666666666666

# Move the active selection to I1 (matches the diff's sheetView selection)
$ws.Range("I1").Select()
